# This workbook is a BFO-2020 relations table. The ISO-release edit trims the
# redundant trailing "at" token from several relation-name cells/shared
# strings (e.g. "has continuant part at" -> "has continuant part "), and
# fixes a typo ("inverse at all times" -> "reverse at all times") in the
# header row. Row order / row count / dimension are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=1; Col="E"; Text="reverse at all times"},
    @{Row=3; Col="A"; Text="concretizes "},
    @{Row=3; Col="D"; Text="is concretized by "},
    @{Row=4; Col="A"; Text="continuant part of "},
    @{Row=4; Col="D"; Text="has continuant part "},
    @{Row=7; Col="A"; Text="generically depends on "},
    @{Row=7; Col="D"; Text="is carrier of "},
    @{Row=8; Col="A"; Text="has continuant part "},
    @{Row=8; Col="D"; Text="continuant part of "},
    @{Row=12; Col="A"; Text="has location "},
    @{Row=12; Col="D"; Text="located in "},
    @{Row=13; Col="A"; Text="has material basis "},
    @{Row=13; Col="D"; Text="material basis of "},
    @{Row=14; Col="A"; Text="has member part "},
    @{Row=14; Col="D"; Text="member part of "},
    @{Row=16; Col="A"; Text="has participant "},
    @{Row=16; Col="D"; Text="participates in "},
    @{Row=17; Col="A"; Text="has proper continuant part "},
    @{Row=17; Col="D"; Text="proper continuant part of "},
    @{Row=24; Col="A"; Text="is carrier of "},
    @{Row=24; Col="D"; Text="generically depends on "},
    @{Row=25; Col="A"; Text="is concretized by "},
    @{Row=25; Col="D"; Text="concretizes "},
    @{Row=27; Col="A"; Text="located in "},
    @{Row=27; Col="D"; Text="has location "},
    @{Row=28; Col="A"; Text="material basis of "},
    @{Row=28; Col="D"; Text="has material basis "},
    @{Row=29; Col="A"; Text="member part of "},
    @{Row=29; Col="D"; Text="has member part "},
    @{Row=30; Col="A"; Text="occupies spatial region "},
    @{Row=35; Col="A"; Text="participates in "},
    @{Row=35; Col="D"; Text="has participant "},
    @{Row=38; Col="A"; Text="proper continuant part of "},
    @{Row=38; Col="D"; Text="has proper continuant part "},
    @{Row=42; Col="A"; Text="spatially projects onto "}
)

foreach ($ch in $changes) {
    $addr = "$($ch.Col)$($ch.Row)"
    $ws.Range($addr).Value = $ch.Text
}

# Move the active selection as it appears in the saved workbook.
$ws.Range("D31").Select()
